# Refresh cryptos list values (coin prices / 1h volume %) per latest scrape.
# Leading apostrophe forces text entry (matches original inlineStr cells)
# without disturbing each cell's NumberFormat/style, exactly like typing
# an apostrophe-prefixed value into Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.765.04"
$ws.Range("E2").Value = "'  -5.42%  "

# Row 3
$ws.Range("D3").Value = "'2.543.91"
$ws.Range("E3").Value = "'  -2.51%  "

# Row 4
$ws.Range("E4").Value = "'  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'301.15"
$ws.Range("E5").Value = "'  -2.01%  "

# Row 6
$ws.Range("D6").Value = "'95.96"
$ws.Range("E6").Value = "'  -4.19%  "

# Row 7
$ws.Range("D7").Value = "'0.579"
$ws.Range("E7").Value = "'  -3.66%  "

# Row 8
$ws.Range("E8").Value = "'  +0.33%  "

# Row 9
$ws.Range("D9").Value = "'0.560"
$ws.Range("E9").Value = "'  -3.00%  "

# Row 10
$ws.Range("D10").Value = "'36.92"
$ws.Range("E10").Value = "'  -6.45%  "

# Row 11
$ws.Range("D11").Value = "'0.0814"
$ws.Range("E11").Value = "'  -3.63%  "

# Row 12
$ws.Range("D12").Value = "'7.78"
$ws.Range("E12").Value = "'  -4.70%  "

# Row 13
$ws.Range("D13").Value = "'2.946.97"
$ws.Range("E13").Value = "'  -1.64%  "

# Row 14
$ws.Range("E14").Value = "'  +1.40%  "

# Row 15
$ws.Range("D15").Value = "'2.563.71"
$ws.Range("E15").Value = "'  -1.33%  "

# Row 16
$ws.Range("D16").Value = "'0.886"
$ws.Range("E16").Value = "'  -3.70%  "

# Row 17
$ws.Range("D17").Value = "'14.29"
$ws.Range("E17").Value = "'  -4.55%  "

# Row 18
$ws.Range("D18").Value = "'43.766.98"
$ws.Range("E18").Value = "'  -5.75%  "

# Row 19
$ws.Range("B19").Value = "'ShibaInu"
$ws.Range("C19").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.0₃0980"
$ws.Range("E19").Value = "'  -3.17%  "

# Row 20
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.66"
$ws.Range("E20").Value = "'  -1.25%  "

# Row 21
$ws.Range("D21").Value = "'12.48"
$ws.Range("E21").Value = "'  -3.55%  "

# Row 22
$ws.Range("D22").Value = "'73.24"
$ws.Range("E22").Value = "'  +2.37%  "

# Row 23
$ws.Range("D23").Value = "'264.14"
$ws.Range("E23").Value = "'  -3.62%  "

# Row 24
$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "'  -3.46%  "

# Row 25
$ws.Range("D25").Value = "'2.19"
$ws.Range("E25").Value = "'  +1.50%  "

# Row 26
$ws.Range("D26").Value = "'29.02"
$ws.Range("E26").Value = "'  +0.81%  "

# Row 27
$ws.Range("E27").Value = "'  +0.05%  "

# Row 28
$ws.Range("D28").Value = "'10.21"
$ws.Range("E28").Value = "'  -3.83%  "

# Row 29
$ws.Range("E29").Value = "'  -2.90%  "

# Row 30
$ws.Range("D30").Value = "'37.95"
$ws.Range("E30").Value = "'  -2.82%  "

# Row 31
$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = "'  -3.68%  "

# Row 32
$ws.Range("D32").Value = "'3.54"
$ws.Range("E32").Value = "'  -2.75%  "

# Row 33
$ws.Range("B33").Value = "'Monero"
$ws.Range("C33").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'152.35"
$ws.Range("E33").Value = "'  +0.93%  "

# Row 34
$ws.Range("B34").Value = "'WEMIXToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.79"
$ws.Range("E34").Value = "'  -2.26%  "

# Row 35
$ws.Range("B35").Value = "'ARBITRUM"
$ws.Range("C35").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'2.18"
$ws.Range("E35").Value = "'  -2.47%  "

# Row 36
$ws.Range("D36").Value = "'0.0809"
$ws.Range("E36").Value = "'  -3.75%  "

# Row 37
$ws.Range("E37").Value = "'  -4.04%  "

# Row 38
$ws.Range("D38").Value = "'0.120"
$ws.Range("E38").Value = "'  -2.40%  "

# Row 39
$ws.Range("D39").Value = "'23.84"
$ws.Range("E39").Value = "'  +1.52%  "

# Row 40
$ws.Range("D40").Value = "'16.64"
$ws.Range("E40").Value = "'  +5.08%  "

# Row 41
$ws.Range("D41").Value = "'3.57"
$ws.Range("E41").Value = "'  -2.02%  "

# Row 42
$ws.Range("D42").Value = "'0.0315"
$ws.Range("E42").Value = "'  -4.85%  "

# Row 43
$ws.Range("D43").Value = "'3.85"
$ws.Range("E43").Value = "'  -5.08%  "

# Row 44
$ws.Range("D44").Value = "'2.024.15"
$ws.Range("E44").Value = "'  -5.06%  "

# Row 45
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "'  -0.17%  "

# Row 46
$ws.Range("D46").Value = "'87.21"
$ws.Range("E46").Value = "'  -6.36%  "

# Row 47
$ws.Range("B47").Value = "'ApeXProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'1.64"
$ws.Range("E47").Value = "'  +5.85%  "

# Row 48
$ws.Range("B48").Value = "'FraxShare"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'9.11"
$ws.Range("E48").Value = "'  -4.15%  "

# Row 49
$ws.Range("D49").Value = "'2.816.27"
$ws.Range("E49").Value = "'  -1.25%  "

# Row 50
$ws.Range("D50").Value = "'104.55"
$ws.Range("E50").Value = "'  -4.11%  "

# Row 51
$ws.Range("D51").Value = "'0.190"
$ws.Range("E51").Value = "'  -4.64%  "

